# Generate Report for Handoff
#
# The "d0c6ac18-ceaf-4b0d-9442-1483750bf9e2" file has dropped out of the
# handoff report, and the "7cfe3a7e-a51c-4ecd-9574-037b5e9db882" file is now
# reported as "Ready for handoff" (rather than "Handed back: in sync with
# en-US"), with refreshed handoff timestamps. Concretely, row 3 (the
# d0c6ac18... entry) is removed from every sheet, and row 2's status/date
# cells are updated on each sheet.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Hyperlinks.Delete()
$ws1.Range("A3:D3").EntireRow.Delete()
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-31-11 12:31:51"
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/263e490a60a3201f12386b49c8b495d1f0e4b396/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md") | Out-Null
$ws1.Range("A2").Style = "HyperLink"

# ---- zh-cn sheet ----
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()
$ws2.Range("A3:K3").EntireRow.Delete()
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("E2").Value = "2016-03-11 12:31:48"
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/263e490a60a3201f12386b49c8b495d1f0e4b396/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/263e490a60a3201f12386b49c8b495d1f0e4b396/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7bbc0d6be9e3c8303a4178fd73f4f0dd932077e6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/83fdb205689054a6fb28f3325f82e5260068bfef/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/65cba7e910ec1c9b11fdbe9cf72b230847e2de32/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.zh-cn.xlf", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.zh-cn.xlf") | Out-Null
$ws2.Range("A2").Style = "HyperLink"
$ws2.Range("B2").Style = "HyperLink"
$ws2.Range("D2").Style = "HyperLink"
$ws2.Range("F2").Style = "HyperLink"
$ws2.Range("G2").Style = "HyperLink"

# ---- de-de sheet ----
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()
$ws3.Range("A3:K3").EntireRow.Delete()
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("E2").Value = "2016-03-11 12:31:51"
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/263e490a60a3201f12386b49c8b495d1f0e4b396/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/263e490a60a3201f12386b49c8b495d1f0e4b396/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7e77109eaa64bd478cdc6319a0d8797ae3b763b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/44b53749e5192a4b812c24a6f90e9466e31dd13d/e2e/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f0cbf690af4b6bc62ed287a5d0a421d7518c1b46/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.de-de.xlf", [System.Type]::Missing, [System.Type]::Missing, "7cfe3a7e-a51c-4ecd-9574-037b5e9db882.b5bd690eba38cdff5d8eb2f78f8401dd2c3ddae5.de-de.xlf") | Out-Null
$ws3.Range("A2").Style = "HyperLink"
$ws3.Range("B2").Style = "HyperLink"
$ws3.Range("D2").Style = "HyperLink"
$ws3.Range("F2").Style = "HyperLink"
$ws3.Range("G2").Style = "HyperLink"
